$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last Refresh Date" header cell
$ws.Range("A1").Value = "Last Refresh Date :08/27/2024"

# Set the Utility Code column (D) to AZ001 for all data rows (4 through 265)
$ws.Range("D4:D265").Value = "AZ001"
